# Update battery costs and endo learn parameters, recalibrate against ATB
#
# Targets the "PDiBCpDoC" (Perc Decline in Battery Cost per Doubling of
# Capacity) workbook:
#   - "About" sheet: refresh the cited source from the old BNEF blog post
#     to BNEF's "Electric Vehicle Outlook 2024", bump the cited year to
#     2024, drop the stray formatted-but-empty cell at D14, and move the
#     selection to reflect the now-shorter content block.
#   - "PDiBCpDoC" sheet: recalibrate the % decline per doubling parameter
#     from 0.18 to 0.17, and update the selection accordingly.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("PDiBCpDoC")

# --- "PDiBCpDoC" data sheet -------------------------------------------------

# Recalibrated % decline per doubling of capacity.
$wsData.Range("B2").Value = 0.17

# Move the selection (no longer parked on the old, unrelated H30 cell).
$wsData.Range("B3").Select()

# --- "About" sheet -----------------------------------------------------------

# Citation year refreshed to 2024.
$wsAbout.Range("B4").Value = 2024

# Swap in the new source title/link (BNEF's Electric Vehicle Outlook 2024),
# replacing the old "Behind the Scenes" blog post citation.
$wsAbout.Range("B5").Value = "Electric Vehicle Outlook 2024"
$wsAbout.Range("B6").Value = "https://about.bnef.com/electric-vehicle-outlook/"

# The stray, empty-but-styled D14 cell is removed entirely (along with its
# now-unused style), shrinking the sheet's used range down to A1:C8.
$wsAbout.Range("D14").Clear()

# Re-activate the About sheet and move the selection to reflect the new,
# shorter layout.
$wsAbout.Activate()
$wsAbout.Range("B6").Select()
